# Updated cryptos list (Price + Volume(1h) columns)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.284.37"
$ws.Range("E2").Value = "  -2.81%  "

$ws.Range("D3").Value = "'3.678.07"
$ws.Range("E3").Value = "  -3.31%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").Value = "'681.72"
$ws.Range("E5").Value = "  -3.29%  "

$ws.Range("D6").Value = "'159.33"
$ws.Range("E6").Value = "  -6.94%  "

$ws.Range("D7").Value = "'3.677.30"
$ws.Range("E7").Value = "  -3.34%  "

$ws.Range("E8").Value = "  -0.06%  "

$ws.Range("D9").Value = "'0.493"
$ws.Range("E9").Value = "  -6.15%  "

$ws.Range("D10").Value = "'0.145"
$ws.Range("E10").Value = "  -9.42%  "

$ws.Range("D11").Value = "'7.14"
$ws.Range("E11").Value = "  -4.11%  "

$ws.Range("D12").Value = "'0.435"
$ws.Range("E12").Value = "  -10.09%  "

$ws.Range("D13").Value = "'0.0000232"
$ws.Range("E13").Value = "  -7.16%  "

$ws.Range("D14").Value = "'4.296.69"
$ws.Range("E14").Value = "  -3.34%  "

$ws.Range("D15").Value = "'32.40"
$ws.Range("E15").Value = "  -10.89%  "

$ws.Range("D16").Value = "'3.678.03"
$ws.Range("E16").Value = "  -2.74%  "

$ws.Range("D17").Value = "'69.319.27"
$ws.Range("E17").Value = "  -2.94%  "

$ws.Range("E18").Value = "  -1.35%  "

$ws.Range("D19").Value = "'15.84"
$ws.Range("E19").Value = "  -9.39%  "

$ws.Range("D20").Value = "'6.41"
$ws.Range("E20").Value = "  -11.07%  "

$ws.Range("D21").Value = "'473.71"
$ws.Range("E21").Value = "  -8.11%  "

$ws.Range("D22").Value = "'9.87"
$ws.Range("E22").Value = "  -5.45%  "

$ws.Range("D23").Value = "'0.648"
$ws.Range("E23").Value = "  -9.25%  "

$ws.Range("D24").Value = "'79.23"
$ws.Range("E24").Value = "  -5.63%  "

$ws.Range("D25").Value = "'3.820.27"
$ws.Range("E25").Value = "  -3.24%  "

$ws.Range("E26").Value = "  +0.02%  "

$ws.Range("D27").Value = "'0.0000124"
$ws.Range("E27").Value = "  -11.92%  "

$ws.Range("D28").Value = "'10.85"
$ws.Range("E28").Value = "  -13.98%  "

$ws.Range("D29").Value = "'9.15"
$ws.Range("E29").Value = "  -11.19%  "

$ws.Range("E30").Value = "  -10.92%  "

$ws.Range("D31").Value = "'1.72"
$ws.Range("E31").Value = "  -14.57%  "

$ws.Range("D32").Value = "'6.64"
$ws.Range("E32").Value = "  -9.99%  "

$ws.Range("D33").Value = "'2.02"
$ws.Range("E33").Value = "  -10.14%  "

$ws.Range("D34").Value = "'0.999"
$ws.Range("E34").Value = "  +0.03%  "

$ws.Range("D35").Value = "'26.61"
$ws.Range("E35").Value = "  -8.71%  "

$ws.Range("D36").Value = "'0.160"
$ws.Range("E36").Value = "  -7.56%  "

$ws.Range("D37").Value = "'8.12"
$ws.Range("E37").Value = "  -12.62%  "

$ws.Range("D38").Value = "'6.04"
$ws.Range("E38").Value = "  -7.53%  "

$ws.Range("E39").Value = "  -0.02%  "

$ws.Range("D40").Value = "'2.23"
$ws.Range("E40").Value = "  -9.13%  "

$ws.Range("D41").Value = "'0.0899"
$ws.Range("E41").Value = "  -11.22%  "

$ws.Range("E42").Value = "  -0.05%  "

$ws.Range("D43").Value = "'0.939"
$ws.Range("E43").Value = "  -7.00%  "

$ws.Range("D44").Value = "'164.87"
$ws.Range("E44").Value = "  -1.74%  "

$ws.Range("D45").Value = "'47.72"
$ws.Range("E45").Value = "  -4.82%  "

$ws.Range("D46").Value = "'2.69"
$ws.Range("E46").Value = "  -16.83%  "

$ws.Range("D47").Value = "'1.30"
$ws.Range("E47").Value = "  -7.03%  "

$ws.Range("D48").Value = "'0.000272"
$ws.Range("E48").Value = "  -11.21%  "

$ws.Range("D49").Value = "'27.96"
$ws.Range("E49").Value = "  -4.79%  "

$ws.Range("D50").Value = "'1.08"
$ws.Range("E50").Value = "  -6.42%  "

$ws.Range("D51").Value = "'7.85"
$ws.Range("E51").Value = "  -9.27%  "
